$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round the Ost/Nord coordinate values to whole numbers
$ws.Range("Q2").Value = 575012
$ws.Range("R2").Value = 6300320

# Remove the Starttid (Z2) and Sluttid (AB2) values entirely
$ws.Range("Z2").ClearContents()
$ws.Range("AB2").ClearContents()
